$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.19925705524686776,
    -0.005999999985213833,
    -0.0039999999868047809,
    -0.0079999999761035667,
    -0.0029999999867262872,
    -0.0019999999862729823,
    -0.0099999999671074313,
    -0.0099999999665776329,
    -0.0019999999851187944,
    -0.0019999999847986061,
    -0.0029999999824337209,
    -0.0034999999812255211,
    -0.0034999999812832527,
    0.013562891188107606,
    -0.00099999998708266702,
    -0.0019999999844833027,
    -0.0019999999841822103,
    -0.0039999999793609575,
    -0.0039999999898872041,
    -0.0039999999892224025,
    -0.0039999999891273674,
    -0.0039999999890385496,
    -0.034042946223591386,
    -0.019999999946414881,
    -0.019999999945747859,
    -0.0024999999847565846,
    -0.0024999999842703069,
    -0.0019999999834077187,
    0.036819312267335036,
    -0.059999999843069141,
    -0.0069999999676522151,
    -0.0099999999603728185,
    -0.003999999974524826
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

